{"js": "// Append four new rows to the end of the (single) table in the document body.\n// Each new row has one cell containing a single line of text, matching the\n// formatting of the table's existing rows (10070 dxa wide column, no wrap,\n// sz/szCs 22 half-points run & paragraph formatting).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRowsText = [\n  [\"CST Warp Engine & Navigation -Galactic Corrections and Core Physics Toolkit\"],\n  [\"Casanova Warp-Drive & CST Navigation -Equations, Laws, and Original Constructs\"],\n  [\"A Theoretical Analysis of the Reported India UFO Event\"],\n  [\"Quantum Credit Energy (QC): An Inclusive, Energy-Backed Civic Economy for Recovery, Reentry, and Community Resilience\"]\n];\n\ntable.addRows(\"End\", newRowsText.length, newRowsText);\nawait context.sync();\n", "ps1": "# Append four new rows to the end of the (single) table in the document,\n# one line of text per row, matching the table's existing row formatting.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRowsText = @(\n  \"CST Warp Engine & Navigation -Galactic Corrections and Core Physics Toolkit\",\n  \"Casanova Warp-Drive & CST Navigation -Equations, Laws, and Original Constructs\",\n  \"A Theoretical Analysis of the Reported India UFO Event\",\n  \"Quantum Credit Energy (QC): An Inclusive, Energy-Backed Civic Economy for Recovery, Reentry, and Community Resilience\"\n)\n\nforeach ($txt in $newRowsText) {\n  $row = $t.Rows.Add()\n  $row.Cells.Item(1).Range.Text = $txt\n}\n"}
